# Generate Report for Handoff
# Updates the status from "In Translation" to "Ready for handoff" and
# refreshes the "Latest HO Xliff Generate Date" / "Latest Handoff Datetime"
# timestamps, then widens the status/date columns to fit the new text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Overview sheet ---------------------------------------------------
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2016-08-13 03:11:52"

$overview.Columns.Item(5).ColumnWidth = 17.2159881591797
$overview.Columns.Item(6).ColumnWidth = 17.2159881591797

# --- zh-cn sheet --------------------------------------------------------
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-08-13 03:11:45"

$zhcn.Columns.Item(3).ColumnWidth = 17.2159881591797

# --- de-de sheet --------------------------------------------------------
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2016-08-13 03:11:52"

$dede.Columns.Item(3).ColumnWidth = 17.2159881591797
